# Update the LR-pairs data with newly recomputed TPM-based statistics.
# Workbook is already open; grab the active worksheet that holds the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Wnt5a -> Fzd4 -> ECs)
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.532132
$ws.Range("N2").Value = 79.596396
$ws.Range("O2").Value = 0.3960736634233649
$ws.Range("P2").Value = 0.3960736634233648
$ws.Range("Q2").Value = 242.953653750632
$ws.Range("R2").Value = 2186.582883755688
$ws.Range("S2").Value = 0.3839425972069762
$ws.Range("T2").Value = 0.3839425972069761

# Row 3 (FAPs -> Wnt5a -> Fzd4 -> FAPs)
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.2505213219764053
$ws.Range("P3").Value = 0.2505213219764053
$ws.Range("S3").Value = 0.2428482777268953
$ws.Range("T3").Value = 0.2428482777268953

# Row 4 (FAPs -> Wnt5a -> Fzd4 -> MuSCs)
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 23.67385
$ws.Range("N4").Value = 71.02154999999999
$ws.Range("O4").Value = 0.3534050146002298
$ws.Range("P4").Value = 0.3534050146002298
$ws.Range("Q4").Value = 216.7804817134333
$ws.Range("R4").Value = 1951.0243354209
$ws.Range("S4").Value = 0.3425808169086589
$ws.Range("T4").Value = 0.3425808169086589

# Row 5 (MuSCs -> Wnt5a -> Fzd4 -> ECs)
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.532132
$ws.Range("N5").Value = 79.596396
$ws.Range("O5").Value = 0.3960736634233649
$ws.Range("P5").Value = 0.3960736634233648
$ws.Range("Q5").Value = 7.676373714724001
$ws.Range("R5").Value = 69.087363432516
$ws.Range("S5").Value = 0.01213106621638873
$ws.Range("T5").Value = 0.01213106621638872

# Row 6 (MuSCs -> Wnt5a -> Fzd4 -> FAPs)
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.2505213219764053
$ws.Range("P6").Value = 0.2505213219764053
$ws.Range("Q6").Value = 4.855398044837889
$ws.Range("R6").Value = 43.698582403541
$ws.Range("S6").Value = 0.007673044249510009
$ws.Range("T6").Value = 0.007673044249510009

# Row 7 (MuSCs -> Wnt5a -> Fzd4 -> MuSCs)
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 23.67385
$ws.Range("N7").Value = 71.02154999999999
$ws.Range("O7").Value = 0.3534050146002298
$ws.Range("P7").Value = 0.3534050146002298
$ws.Range("Q7").Value = 6.849405086116667
$ws.Range("R7").Value = 61.64464577504999
$ws.Range("S7").Value = 0.01082419769157089
$ws.Range("T7").Value = 0.01082419769157089
